$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
$ws.Range("B1:J1").EntireColumn.Delete()
$ws.Range("A1").ClearContents()
